$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(34, 8).Value = 7598   # H34
$ws.Cells.Item(34, 9).Value = 7969.7144   # I34
$ws.Cells.Item(34, 10).Value = 4996   # J34
$ws.Cells.Item(34, 11).Value = 7969.7144   # K34
$ws.Cells.Item(34, 12).Value = 4996   # L34
$ws.Cells.Item(34, 13).Value = -7766.7144   # M34
$ws.Cells.Item(34, 14).Value = -5402   # N34
$ws.Cells.Item(36, 8).Value = 7598   # H36
$ws.Cells.Item(36, 9).Value = 7969.7144   # I36
$ws.Cells.Item(36, 10).Value = 4996   # J36
$ws.Cells.Item(36, 11).Value = 7969.7144   # K36
$ws.Cells.Item(36, 12).Value = 4996   # L36
$ws.Cells.Item(36, 13).Value = -7254.7144   # M36
$ws.Cells.Item(36, 14).Value = -6426   # N36
$ws.Cells.Item(74, 8).Value = 8425.25   # H74
$ws.Cells.Item(74, 9).Value = 5800.6665   # I74
$ws.Cells.Item(74, 11).Value = 5800.6665   # K74
$ws.Cells.Item(74, 13).Value = -4864.6665   # M74
$ws.Cells.Item(77, 8).Value = 8425.25   # H77
$ws.Cells.Item(77, 9).Value = 5800.6665   # I77
$ws.Cells.Item(77, 11).Value = 29003.3325   # K77
$ws.Cells.Item(77, 13).Value = -24323.3325   # M77
$ws.Cells.Item(100, 8).Value = 5999.9   # H100
$ws.Cells.Item(100, 9).Value = 1066.6666   # I100
$ws.Cells.Item(100, 10).Value = 8114.143   # J100
$ws.Cells.Item(100, 11).Value = 1066.6666   # K100
$ws.Cells.Item(100, 12).Value = 8114.143   # L100
$ws.Cells.Item(100, 13).Value = -525.6666   # M100
$ws.Cells.Item(100, 14).Value = -9196.143   # N100
$ws.Cells.Item(103, 8).Value = 1546.1538   # H103
$ws.Cells.Item(103, 9).Value = 566.6667   # I103
$ws.Cells.Item(103, 10).Value = 1840   # J103
$ws.Cells.Item(103, 11).Value = 1700.0001   # K103
$ws.Cells.Item(103, 12).Value = 5520   # L103
$ws.Cells.Item(103, 13).Value = -1114.0001   # M103
$ws.Cells.Item(103, 14).Value = -6692   # N103
$ws.Cells.Item(129, 8).Value = 2206.5293   # H129
$ws.Cells.Item(129, 9).Value = 1543   # I129
$ws.Cells.Item(129, 11).Value = 4629   # K129
$ws.Cells.Item(129, 13).Value = 371   # M129
$ws.Cells.Item(138, 8).Value = 4489.5654   # H138
$ws.Cells.Item(138, 10).Value = 4838.4614   # J138
$ws.Cells.Item(138, 12).Value = 14515.3842   # L138
$ws.Cells.Item(138, 14).Value = -24795.3842   # N138

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 19894.857   # H32
$ws.Cells.Item(32, 9).Value = 19563.59   # I32
$ws.Cells.Item(32, 11).Value = 19563.59   # K32
$ws.Cells.Item(32, 13).Value = -19276.59   # M32
$ws.Cells.Item(110, 8).Value = 10871111   # H110
$ws.Cells.Item(110, 9).Value = 16667370   # I110
$ws.Cells.Item(110, 11).Value = 16667370   # K110
$ws.Cells.Item(110, 13).Value = -16665325   # M110
$ws.Cells.Item(122, 8).Value = 2411.0908   # H122
$ws.Cells.Item(122, 9).Value = 1459.8572   # I122
$ws.Cells.Item(122, 10).Value = 4075.75   # J122
$ws.Cells.Item(122, 11).Value = 4379.571599999999   # K122
$ws.Cells.Item(122, 12).Value = 12227.25   # L122
$ws.Cells.Item(122, 13).Value = -1929.571599999999   # M122
$ws.Cells.Item(122, 14).Value = -17127.25   # N122
$ws.Cells.Item(132, 8).Value = 3871.0393   # H132
$ws.Cells.Item(132, 9).Value = 3386.1836   # I132
$ws.Cells.Item(132, 11).Value = 10158.5508   # K132
$ws.Cells.Item(132, 13).Value = -7628.550799999999   # M132

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 112478.945   # H86
$ws.Cells.Item(86, 9).Value = 1594.1538   # I86
$ws.Cells.Item(86, 11).Value = 1594.1538   # K86
$ws.Cells.Item(86, 13).Value = -471.1538   # M86
$ws.Cells.Item(89, 8).Value = 112478.945   # H89
$ws.Cells.Item(89, 9).Value = 1594.1538   # I89
$ws.Cells.Item(89, 11).Value = 7970.769   # K89
$ws.Cells.Item(89, 13).Value = -2354.769   # M89
$ws.Cells.Item(99, 8).Value = 2167.2942   # H99
$ws.Cells.Item(99, 9).Value = 1793.7778   # I99
$ws.Cells.Item(99, 11).Value = 1793.7778   # K99
$ws.Cells.Item(99, 13).Value = -295.7778000000001   # M99
$ws.Cells.Item(105, 8).Value = 35726532   # H105
$ws.Cells.Item(105, 9).Value = 62518276   # I105
$ws.Cells.Item(105, 11).Value = 62518276   # K105
$ws.Cells.Item(105, 13).Value = -62516529   # M105
$ws.Cells.Item(137, 8).Value = 120000   # H137
$ws.Cells.Item(137, 10).Value = 120000   # J137
$ws.Cells.Item(137, 12).Value = 120000   # L137
$ws.Cells.Item(137, 14).Value = -130200   # N137

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(37, 8).Value = 68298.92   # H37
$ws.Cells.Item(37, 10).Value = 68298.92   # J37
$ws.Cells.Item(37, 12).Value = 204896.76   # L37
$ws.Cells.Item(37, 14).Value = -205120.76   # N37
$ws.Cells.Item(54, 8).Value = 3187.5   # H54
$ws.Cells.Item(54, 9).Value = 750   # I54
$ws.Cells.Item(54, 11).Value = 2250   # K54
$ws.Cells.Item(54, 13).Value = -1691   # M54
$ws.Cells.Item(62, 8).Value = 1892.3036   # H62
$ws.Cells.Item(62, 9).Value = 999   # I62
$ws.Cells.Item(62, 11).Value = 2997   # K62
$ws.Cells.Item(62, 13).Value = -2311   # M62
$ws.Cells.Item(65, 8).Value = 1892.3036   # H65
$ws.Cells.Item(65, 9).Value = 999   # I65
$ws.Cells.Item(65, 11).Value = 8991   # K65
$ws.Cells.Item(65, 13).Value = -5559   # M65
$ws.Cells.Item(87, 8).Value = 199.5   # H87
$ws.Cells.Item(87, 9).Value = 199.5   # I87
$ws.Cells.Item(87, 11).Value = 598.5   # K87
$ws.Cells.Item(87, 13).Value = 649.5   # M87
$ws.Cells.Item(90, 8).Value = 199.5   # H90
$ws.Cells.Item(90, 9).Value = 199.5   # I90
$ws.Cells.Item(90, 11).Value = 1795.5   # K90
$ws.Cells.Item(90, 13).Value = 4444.5   # M90
$ws.Cells.Item(140, 8).Value = 17975.666   # H140
$ws.Cells.Item(140, 9).Value = 38437   # I140
$ws.Cells.Item(140, 10).Value = 1606.6   # J140
$ws.Cells.Item(140, 11).Value = 115311   # K140
$ws.Cells.Item(140, 12).Value = 4819.799999999999   # L140
$ws.Cells.Item(140, 13).Value = -110131   # M140
$ws.Cells.Item(140, 14).Value = -15179.8   # N140

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(36, 8).Value = 5347.143   # H36
$ws.Cells.Item(36, 9).Value = 2725   # I36
$ws.Cells.Item(36, 10).Value = 6396   # J36
$ws.Cells.Item(36, 11).Value = 2725   # K36
$ws.Cells.Item(36, 12).Value = 6396   # L36
$ws.Cells.Item(36, 13).Value = -2240   # M36
$ws.Cells.Item(36, 14).Value = -7366   # N36
$ws.Cells.Item(70, 8).Value = 5885.2856   # H70
$ws.Cells.Item(70, 9).Value = 5502.6665   # I70
$ws.Cells.Item(70, 11).Value = 5502.6665   # K70
$ws.Cells.Item(70, 13).Value = -5232.6665   # M70
$ws.Cells.Item(73, 8).Value = 5885.2856   # H73
$ws.Cells.Item(73, 9).Value = 5502.6665   # I73
$ws.Cells.Item(73, 11).Value = 5502.6665   # K73
$ws.Cells.Item(73, 13).Value = -4566.6665   # M73
$ws.Cells.Item(122, 8).Value = 2327.8   # H122
$ws.Cells.Item(122, 10).Value = 2596.5334   # J122
$ws.Cells.Item(122, 12).Value = 7789.600199999999   # L122
$ws.Cells.Item(122, 14).Value = -12689.6002   # N122
$ws.Cells.Item(132, 8).Value = 5209.8486   # H132
$ws.Cells.Item(132, 9).Value = 4097.6   # I132
$ws.Cells.Item(132, 10).Value = 16332.333   # J132
$ws.Cells.Item(132, 11).Value = 12292.8   # K132
$ws.Cells.Item(132, 12).Value = 48996.999   # L132
$ws.Cells.Item(132, 13).Value = -9762.800000000001   # M132
$ws.Cells.Item(132, 14).Value = -54056.999   # N132

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 868.35297   # H16
$ws.Cells.Item(16, 10).Value = 449.66666   # J16
$ws.Cells.Item(16, 12).Value = 449.66666   # L16
$ws.Cells.Item(16, 14).Value = -789.66666   # N16
$ws.Cells.Item(40, 8).Value = 19238464   # H40
$ws.Cells.Item(40, 10).Value = 10480.4   # J40
$ws.Cells.Item(40, 12).Value = 10480.4   # L40
$ws.Cells.Item(40, 14).Value = -10752.4   # N40
$ws.Cells.Item(82, 8).Value = 961.7143   # H82
$ws.Cells.Item(82, 10).Value = 697.36365   # J82
$ws.Cells.Item(82, 12).Value = 697.36365   # L82
$ws.Cells.Item(82, 14).Value = -1419.36365   # N82
$ws.Cells.Item(85, 8).Value = 961.7143   # H85
$ws.Cells.Item(85, 10).Value = 697.36365   # J85
$ws.Cells.Item(85, 12).Value = 697.36365   # L85
$ws.Cells.Item(85, 14).Value = -3193.36365   # N85
$ws.Cells.Item(93, 8).Value = 2132.4   # H93
$ws.Cells.Item(93, 10).Value = 2923.7   # J93
$ws.Cells.Item(93, 12).Value = 2923.7   # L93
$ws.Cells.Item(93, 14).Value = -5419.7   # N93
$ws.Cells.Item(97, 8).Value = 0   # H97
$ws.Cells.Item(97, 10).Value = 0   # J97
$ws.Cells.Item(97, 14).ClearContents()   # N97 (removed)
$ws.Cells.Item(100, 8).Value = 8931584   # H100
$ws.Cells.Item(100, 9).Value = 50001600   # I100
$ws.Cells.Item(100, 11).Value = 50001600   # K100
$ws.Cells.Item(100, 13).Value = -50001059   # M100
$ws.Cells.Item(122, 8).Value = 4318   # H122
$ws.Cells.Item(122, 9).Value = 4318   # I122
$ws.Cells.Item(122, 11).Value = 12954   # K122
$ws.Cells.Item(122, 13).Value = -10504   # M122
$ws.Cells.Item(136, 8).Value = 5199.4   # H136
$ws.Cells.Item(136, 9).Value = 3991.6155   # I136
$ws.Cells.Item(136, 11).Value = 11974.8465   # K136
$ws.Cells.Item(136, 13).Value = -9424.8465   # M136

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(76, 8).Value = 55999.8   # H76
$ws.Cells.Item(76, 10).Value = 55999.8   # J76
$ws.Cells.Item(76, 12).Value = 55999.8   # L76
$ws.Cells.Item(76, 14).Value = -56629.8   # N76
$ws.Cells.Item(79, 8).Value = 55999.8   # H79
$ws.Cells.Item(79, 10).Value = 55999.8   # J79
$ws.Cells.Item(79, 12).Value = 55999.8   # L79
$ws.Cells.Item(79, 14).Value = -58183.8   # N79
$ws.Cells.Item(81, 8).Value = 12334.755   # H81
$ws.Cells.Item(81, 9).Value = 2649.7334   # I81
$ws.Cells.Item(81, 10).Value = 16157.789   # J81
$ws.Cells.Item(81, 11).Value = 5299.4668   # K81
$ws.Cells.Item(81, 12).Value = 32315.578   # L81
$ws.Cells.Item(81, 13).Value = -4238.4668   # M81
$ws.Cells.Item(81, 14).Value = -34437.578   # N81
$ws.Cells.Item(84, 8).Value = 12334.755   # H84
$ws.Cells.Item(84, 9).Value = 2649.7334   # I84
$ws.Cells.Item(84, 10).Value = 16157.789   # J84
$ws.Cells.Item(84, 11).Value = 26497.334   # K84
$ws.Cells.Item(84, 12).Value = 161577.89   # L84
$ws.Cells.Item(84, 13).Value = -21193.334   # M84
$ws.Cells.Item(84, 14).Value = -172185.89   # N84
$ws.Cells.Item(96, 8).Value = 1382.5   # H96
$ws.Cells.Item(96, 9).Value = 1188.8572   # I96
$ws.Cells.Item(96, 11).Value = 1188.8572   # K96
$ws.Cells.Item(96, 13).Value = 184.1428000000001   # M96
